$wb = $excel.ActiveWorkbook

$wsPatient = $wb.Worksheets.Item("patient")
$wsSample  = $wb.Worksheets.Item("patient_sample")
$wsCell    = $wb.Worksheets.Item("cell_model")
$wsShare   = $wb.Worksheets.Item("sharing")

# --- Fix up the model-id values on "cell_model" and "sharing" so they
# --- match the canonical CRC0228PR / CRC0228PRaS model ids used by the
# --- rest of the workbook (patient / pdx_model / model_validation),
# --- instead of the old CRC0014LM* placeholders.

# cell_model!A2/A3 currently carry the "flagged" style (s=24); copy the
# plain text style (s=4) already used by sharing!A2/A3 before overwriting
# the values so the highlighted formatting goes away along with the fix.
$wsShare.Range("A2").Copy()
$wsCell.Range("A2").PasteSpecial(-4122)
$wsShare.Range("A3").Copy()
$wsCell.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsCell.Range("A2").Value = "CRC0228PR"
$wsCell.Range("A3").Value = "CRC0228PRaS"

$wsShare.Range("A2").Value = "CRC0228PR"
$wsShare.Range("A3").Value = "CRC0228PRaS"

# --- Update the saved selections on each sheet.
$wsSample.Activate()
$wsSample.Range("T2:T3").Select()

$wsShare.Activate()
$wsShare.Range("A2:A3").Select()

$wsPatient.Activate()
$wsPatient.Range("H2").Select()

$wsCell.Activate()
$wsCell.Range("A2:A3").Select()
